$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp label (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 13:22"

# --- Swap Sevilla / Gipuzkoa-Guipuzcoa rows (table re-sorted by updated totals) ---
# Row 21 now shows Gipuzkoa/Guipuzcoa's (updated) figures
$ws.Range("A21").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B21").Value = 2086
$ws.Range("C21").Value = 5193
$ws.Range("D21").Value = 5174
$ws.Range("E21").Value = 136

# Row 22 now shows Sevilla's figures
$ws.Range("A22").Value = "Sevilla"
$ws.Range("B22").Value = 2083
$ws.Range("C22").Value = 294
$ws.Range("D22").Value = 1611
$ws.Range("E22").Value = 178

# --- Update Bizkaia/Vizcaya totals (row 7) ---
$ws.Range("B7").Value = 6085
$ws.Range("C7").Value = 5193
$ws.Range("D7").Value = 5174
$ws.Range("E7").Value = 450

# --- Update Araba/Alava totals (row 16) ---
$ws.Range("B16").Value = 3055
$ws.Range("C16").Value = 5193
$ws.Range("D16").Value = 5174
$ws.Range("E16").Value = 273
